$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 155164
$ws.Range("J3").Value = 155164
$ws.Range("L3").Value = 155164
$ws.Range("N3").Value = -155392
$ws.Range("H15").Value = 140231.47
$ws.Range("I15").Value = 140231.47
$ws.Range("K15").Value = 420694.41
$ws.Range("M15").Value = -420525.41
$ws.Range("H74").Value = 4054.9092
$ws.Range("I74").Value = 3800
$ws.Range("J74").Value = 4111.5557
$ws.Range("L74").Value = 4111.5557
$ws.Range("M74").Value = -2864
$ws.Range("N74").Value = -5983.5557
$ws.Range("H76").Value = 37039636
$ws.Range("I76").Value = 111111110
$ws.Range("K76").Value = 111111110
$ws.Range("M76").Value = -111110795
$ws.Range("H77").Value = 4054.9092
$ws.Range("I77").Value = 3800
$ws.Range("J77").Value = 4111.5557
$ws.Range("K77").Value = 19000
$ws.Range("L77").Value = 20557.7785
$ws.Range("M77").Value = -14320
$ws.Range("N77").Value = -29917.7785
$ws.Range("H79").Value = 37039636
$ws.Range("I79").Value = 111111110
$ws.Range("K79").Value = 111111110
$ws.Range("M79").Value = -111110018
$ws.Range("H102").Value = 155164
$ws.Range("J102").Value = 155164
$ws.Range("L102").Value = 155164
$ws.Range("N102").Value = -161654
$ws.Range("H113").Value = 89218.336
$ws.Range("I113").Value = 89218.336
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 89218.336
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -85964.336
$ws.Range("N113").ClearContents()
$ws.Range("H135").Value = 1610.24
$ws.Range("I135").Value = 1586
$ws.Range("K135").Value = 14274
$ws.Range("M135").Value = -11739
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2579.6924
$ws.Range("I61").Value = 1693.375
$ws.Range("K61").Value = 1693.375
$ws.Range("M61").Value = -1481.375
$ws.Range("H88").Value = 11000
$ws.Range("I88").Value = 2000
$ws.Range("J88").Value = 20000
$ws.Range("K88").Value = 2000
$ws.Range("L88").Value = 20000
$ws.Range("M88").Value = -1594
$ws.Range("N88").Value = -20812
$ws.Range("H91").Value = 11000
$ws.Range("I91").Value = 2000
$ws.Range("J91").Value = 20000
$ws.Range("K91").Value = 2000
$ws.Range("L91").Value = 20000
$ws.Range("M91").Value = -596
$ws.Range("N91").Value = -22808
$ws.Range("H136").Value = 2579.6924
$ws.Range("I136").Value = 1693.375
$ws.Range("K136").Value = 5080.125
$ws.Range("M136").Value = -2530.125
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H50").Value = 25350
$ws.Range("J50").Value = 25350
$ws.Range("L50").Value = 25350
$ws.Range("N50").Value = -26498
$ws.Range("H86").Value = 5205
$ws.Range("I86").Value = 1891.7142
$ws.Range("K86").Value = 1891.7142
$ws.Range("M86").Value = -768.7141999999999
$ws.Range("H89").Value = 5205
$ws.Range("I89").Value = 1891.7142
$ws.Range("K89").Value = 9458.571
$ws.Range("M89").Value = -3842.571
$ws.Range("H112").Value = 30000
$ws.Range("J112").Value = 30000
$ws.Range("L112").Value = 30000
$ws.Range("N112").Value = -32954
$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1466.9318
$ws.Range("I31").Value = 957.6667
$ws.Range("J31").Value = 2275.7646
$ws.Range("K31").Value = 957.6667
$ws.Range("L31").Value = 2275.7646
$ws.Range("M31").Value = -662.6667
$ws.Range("N31").Value = -2865.7646
$ws.Range("H34").Value = 1466.9318
$ws.Range("I34").Value = 957.6667
$ws.Range("J34").Value = 2275.7646
$ws.Range("K34").Value = 957.6667
$ws.Range("L34").Value = 2275.7646
$ws.Range("M34").Value = -755.6667
$ws.Range("N34").Value = -2679.7646
$ws.Range("H132").Value = 3300.92
$ws.Range("I132").Value = 2813.375
$ws.Range("K132").Value = 8440.125
$ws.Range("M132").Value = -5910.125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1599
$ws.Range("I5").Value = 1429.7693
$ws.Range("J5").Value = 2039
$ws.Range("K5").Value = 4289.3079
$ws.Range("L5").Value = 6117
$ws.Range("M5").Value = -4177.3079
$ws.Range("N5").Value = -6341
$ws.Range("H118").Value = 2300
$ws.Range("I118").Value = 800
$ws.Range("J118").Value = 3050
$ws.Range("K118").Value = 2400
$ws.Range("L118").Value = 9150
$ws.Range("M118").Value = -1157
$ws.Range("N118").Value = -11636
$ws.Range("H135").Value = 1599
$ws.Range("I135").Value = 1429.7693
$ws.Range("J135").Value = 2039
$ws.Range("K135").Value = 12867.9237
$ws.Range("L135").Value = 18351
$ws.Range("M135").Value = -10332.9237
$ws.Range("N135").Value = -23421
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6598.091
$ws.Range("I70").Value = 6909.875
$ws.Range("J70").Value = 5766.6665
$ws.Range("K70").Value = 6909.875
$ws.Range("L70").Value = 5766.6665
$ws.Range("M70").Value = -6639.875
$ws.Range("N70").Value = -6306.6665
$ws.Range("H73").Value = 6598.091
$ws.Range("I73").Value = 6909.875
$ws.Range("J73").Value = 5766.6665
$ws.Range("K73").Value = 6909.875
$ws.Range("L73").Value = 5766.6665
$ws.Range("M73").Value = -5973.875
$ws.Range("N73").Value = -7638.6665
$ws.Range("H93").Value = 28823.143
$ws.Range("J93").Value = 28823.143
$ws.Range("L93").Value = 28823.143
$ws.Range("N93").Value = -32567.143
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 933.9375
$ws.Range("I46").Value = 970.7778
$ws.Range("J46").Value = 886.5714
$ws.Range("K46").Value = 970.7778
$ws.Range("L46").Value = 886.5714
$ws.Range("M46").Value = -782.7778
$ws.Range("N46").Value = -1262.5714
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 706.5833
$ws.Range("I113").Value = 414.6
$ws.Range("J113").Value = 2166.5
$ws.Range("K113").Value = 1243.8
$ws.Range("L113").Value = 6499.5
$ws.Range("M113").Value = 926.1999999999998
$ws.Range("N113").Value = -10839.5
